# Updated symbol list on Fri Dec 16 17:49:15 UTC 2022 with GitHub Actions
#
# The "Price" column (D) and some "Volume(1h)" column (E) cells in this
# workbook are stored as TEXT even though many of them look like plain
# numbers (t="inlineStr" in the original OOXML). Writing a numeric-looking
# string straight into Range.Value would make Excel silently reinterpret it
# as a number, which would flip the cell's stored type and break the
# text-preserving diff we need to reproduce. To keep the cell text-typed we
# prefix numeric-looking values with a leading apostrophe (the classic
# "force text" trick) and then reset the cell style back to "Normal" so we
# don't leave a stray quote-prefixed style behind on the cell itself.

function Set-TextValue {
    param(
        $Sheet,
        [string]$CellAddr,
        [string]$NewValue
    )

    $range = $Sheet.Range($CellAddr)

    if ($NewValue -match '^[+-]?\d+(\.\d+)?$') {
        # Numeric-looking text: force Excel to keep it as a string, then
        # restore the default "Normal" style so no quotePrefix formatting
        # lingers on the cell itself.
        $range.Value = "'" + $NewValue
        $range.Style = "Normal"
    }
    else {
        # Plain text: no special handling required.
        $range.Value = $NewValue
    }
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

Set-TextValue $ws "D2"  '243.45'
Set-TextValue $ws "D3"  '23.23'
Set-TextValue $ws "D4"  '5.585'
Set-TextValue $ws "D5"  '0.05862'
Set-TextValue $ws "D6"  '3.409'
Set-TextValue $ws "D7"  '6.476'
Set-TextValue $ws "D9"  '0.7983'
Set-TextValue $ws "D10" '0.1461'
Set-TextValue $ws "D11" '0.07633'
Set-TextValue $ws "D12" '0.03245'
Set-TextValue $ws "D13" '0.03000'
Set-TextValue $ws "D14" '0.09242'
Set-TextValue $ws "D15" '0.001666'
Set-TextValue $ws "D16" '3.410'
Set-TextValue $ws "D17" '0.04737'
Set-TextValue $ws "D18" '0.01250'
Set-TextValue $ws "E18" '17OneONEBestin24h'
Set-TextValue $ws "D19" '0.006239'
Set-TextValue $ws "D21" '0.003827'
Set-TextValue $ws "D22" '0.0001502'
Set-TextValue $ws "D23" '3.697'
Set-TextValue $ws "D24" '2.210'
Set-TextValue $ws "D25" '0.3336'
Set-TextValue $ws "D26" '0.1251'
Set-TextValue $ws "D27" '0.0004005'
Set-TextValue $ws "D40" '0.04317'
Set-TextValue $ws "D41" '0.007094'
Set-TextValue $ws "D42" '0.1051'
Set-TextValue $ws "D43" '0.003390'
Set-TextValue $ws "D44" '0.008789'
Set-TextValue $ws "D46" '0.00005753'
Set-TextValue $ws "D48" '0.7862'
Set-TextValue $ws "D49" '0.1048'
Set-TextValue $ws "E49" '48BOLOBOLO'
Set-TextValue $ws "D50" '0.00002103'
